# Append the next batch of cpu_windows hourly CPU-usage log entries
# (rows 148-201) to Sheet1, continuing the existing date_hour_of_Aug /
# cpu_used series. A new label "10_00" that didn't previously exist in
# the shared-string table is introduced at row 161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Label = "6_18"; Value = 70.23 },
    @{ Label = "7_00"; Value = 80.98 },
    @{ Label = "7_06"; Value = 76.51000000000001 },
    @{ Label = "7_12"; Value = 75.67 },
    @{ Label = "7_18"; Value = 78.65000000000001 },
    @{ Label = "8_00"; Value = 77.34999999999999 },
    @{ Label = "8_06"; Value = 40.38 },
    @{ Label = "8_12"; Value = 79.79000000000001 },
    @{ Label = "8_18"; Value = 40.44 },
    @{ Label = "9_00"; Value = 39.22 },
    @{ Label = "9_06"; Value = 90.31 },
    @{ Label = "9_12"; Value = 76.29000000000001 },
    @{ Label = "9_18"; Value = 31.92 },
    @{ Label = "10_00"; Value = 40.5 },
    @{ Label = "10_06"; Value = 55.07 },
    @{ Label = "10_12"; Value = 39.44 },
    @{ Label = "10_18"; Value = 42.75 },
    @{ Label = "11_00"; Value = 35.79 },
    @{ Label = "11_06"; Value = 42.58 },
    @{ Label = "11_12"; Value = 5.05 },
    @{ Label = "11_18"; Value = 4.44 },
    @{ Label = "12_00"; Value = 4.25 },
    @{ Label = "12_06"; Value = 8.07 },
    @{ Label = "12_12"; Value = 7.44 },
    @{ Label = "12_18"; Value = 7.93 },
    @{ Label = "13_00"; Value = 7.79 },
    @{ Label = "13_06"; Value = 7.31 },
    @{ Label = "13_12"; Value = 7.28 },
    @{ Label = "13_18"; Value = 7.27 },
    @{ Label = "14_00"; Value = 7.4 },
    @{ Label = "14_12"; Value = 1.09 },
    @{ Label = "14_18"; Value = 1.69 },
    @{ Label = "15_00"; Value = 1.91 },
    @{ Label = "15_06"; Value = 1.72 },
    @{ Label = "15_12"; Value = 81.76000000000001 },
    @{ Label = "15_18"; Value = 32.15 },
    @{ Label = "16_00"; Value = 41.59 },
    @{ Label = "16_06"; Value = 86.44 },
    @{ Label = "16_12"; Value = 42.14 },
    @{ Label = "16_18"; Value = 33.3 },
    @{ Label = "17_00"; Value = 42.17 },
    @{ Label = "17_06"; Value = 82.48999999999999 },
    @{ Label = "17_12"; Value = 70.86 },
    @{ Label = "17_18"; Value = 81.95999999999999 },
    @{ Label = "18_00"; Value = 86.34 },
    @{ Label = "18_06"; Value = 82.08 },
    @{ Label = "18_12"; Value = 73.84999999999999 },
    @{ Label = "18_18"; Value = 81.89 },
    @{ Label = "19_00"; Value = 84.13 },
    @{ Label = "19_06"; Value = 83.11 },
    @{ Label = "19_12"; Value = 83.16 },
    @{ Label = "19_18"; Value = 81.7 },
    @{ Label = "20_00"; Value = 73.90000000000001 },
    @{ Label = "20_06"; Value = 73.59999999999999 }
)

$startRow = 148
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i].Label
    $ws.Cells.Item($row, 2).Value = $newRows[$i].Value
}

Write-Host "Appended $($newRows.Count) rows (cpu_windows log) -> new range A1:B$($startRow + $newRows.Count - 1)"
